$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the license value in B2 with the latest data
$ws.Range("B2").Value = 59942

# Move / update the active selection to E7 (as last interacted cell)
$ws.Range("E7").Select()
